$wb = $excel.ActiveWorkbook

# --- Sheet: Forecast Comparison ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$ws1.Range("L2").Value = 0.8100000000000001
$ws1.Range("L3").Value = 0.99
$ws1.Range("L4").Value = 1.13
$ws1.Range("L5").Value = 0.82
$ws1.Range("L6").Value = 0.89
$ws1.Range("L7").Value = 0.89
$ws1.Range("L8").Value = 1.02
$ws1.Range("L9").Value = 1.16
$ws1.Range("L10").Value = 1.2
$ws1.Range("L11").Value = 1.15
$ws1.Range("L12").Value = 1.15
$ws1.Range("L13").Value = 1.02
$ws1.Range("L14").Value = 1.05
$ws1.Range("D15").Value = 95
$ws1.Range("L15").Value = 0.86
$ws1.Range("L16").Value = 1.19
$ws1.Range("L17").Value = 0.84

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
# Leading apostrophe keeps this numeric-looking value stored as text,
# matching the other cells in this column (e.g. "1773" before it).
$ws2.Range("B9").Value = "'1772"
